# Update "想去人数" (interested-count) figures on the 展览 and 全部类型 sheets
# to reflect the latest generated output (commit 456a3b4).

$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions) - rows 2,5,6,7,10,11,16,17
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 1206
$ws1.Range("F5").Value = 39
$ws1.Range("F6").Value = 214
$ws1.Range("F7").Value = 67
$ws1.Range("F10").Value = 5624
$ws1.Range("F11").Value = 4988
$ws1.Range("F16").Value = 206
$ws1.Range("F17").Value = 10

# Sheet "全部类型" (all types) - rows 2,5,6,7,10,11,16,19
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 1206
$ws4.Range("F5").Value = 39
$ws4.Range("F6").Value = 214
$ws4.Range("F7").Value = 67
$ws4.Range("F10").Value = 5624
$ws4.Range("F11").Value = 4988
$ws4.Range("F16").Value = 206
$ws4.Range("F19").Value = 10
